$d = $word.ActiveDocument

# The shape with the group of VML shapes lives in the 2nd paragraph of the
# document (the one whose run holds the <w:pict> with the v:group canvas).
$p = $d.Paragraphs(2).Range

# Pull this paragraph's canonical OOXML (round-tripped through a small
# single-paragraph "package") so we can surgically edit the VML inside it.
$full = $p.WordOpenXML

# Isolate the outer <w:p>...</w:p> that is the direct child of <w:body>.
# Nested <w:p> elements exist inside each shape's <w:txbxContent>, so we
# can't just grab the first "</w:p>"; instead we look for the closing
# </w:p> that follows the very last </w:txbxContent> in the payload - that
# is guaranteed to be the outer paragraph's own end tag.
$bodyTag = "<w:body>"
$startIdx = $full.IndexOf($bodyTag) + $bodyTag.Length

$lastTxbx = $full.LastIndexOf("</w:txbxContent>")
$closeTag = "</w:p>"
$closeIdx = $full.IndexOf($closeTag, $lastTxbx)
$endIdx = $closeIdx + $closeTag.Length

$paragraphXml = $full.Substring($startIdx, $endIdx - $startIdx)

# The new third rectangle shape, inserted right after the "Second shape"
# v:shape and before the group's w10:wrap, matching the WW8 import's
# handling of shape text that doesn't fit on a single line.
$newShape = '<v:rect id="_x0000_s1042" style="position:absolute;left:6823;top:1412;width:2094;height:495"><v:textbox><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Third shape with automatically wrapped text.</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect>'

$wrapTag = '<w10:wrap type="none"/>'
$wrapIdx = $paragraphXml.IndexOf($wrapTag)
if ($wrapIdx -lt 0) {
    throw "could not locate w10:wrap anchor inside the group shape paragraph"
}
$updatedParagraphXml = $paragraphXml.Substring(0, $wrapIdx) + $newShape + $paragraphXml.Substring($wrapIdx)

# InsertXML replaces the exact range's contents, so this swaps the whole
# paragraph (pict and all) for the version carrying the extra rectangle.
$result = $p.InsertXML($updatedParagraphXml)

Write-Output "done"
